$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.016.32'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '1.667.65'
$ws.Range('E3').Value = '  -1.54%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.91'
$ws.Range('E5').Value = '  -1.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5107'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2657'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06408'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.82'
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07448'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.671.00'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.511'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5842'
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008582'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '26.095.82'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.78'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.04'
$ws.Range('E21').Value = '  +2.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.211'
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.78'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.623'
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('E26').Value = '  +2.91%  '
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06500'
$ws.Range('E28').Value = '  +13.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.339'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.316'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.518'
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6120'
$ws.Range('E35').Value = '  +2.05%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.681'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.273'
$ws.Range('E38').Value = '  +7.57%  '
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = '1.090.21'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8632'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.92'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').Value = '1.816.40'
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('E45').Value = '  -5.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.49'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.008'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.084'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05233'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4285'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.072'
$ws.Range('E51').Value = '  +4.81%  '
